$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2790.0833
$ws.Range("J17").Value = 2790.0833
$ws.Range("L17").Value = 8370.249899999999
$ws.Range("N17").Value = -8706.249899999999
$ws.Range("H62").Value = 1000
$ws.Range("I62").Value = 1000
$ws.Range("K62").Value = 1000
$ws.Range("M62").Value = -376
$ws.Range("H65").Value = 1000
$ws.Range("I65").Value = 1000
$ws.Range("K65").Value = 5000
$ws.Range("M65").Value = -1880
$ws.Range("H80").Value = 284.5
$ws.Range("J80").Value = 246.66667
$ws.Range("L80").Value = 740.00001
$ws.Range("N80").Value = -2736.00001
$ws.Range("H83").Value = 284.5
$ws.Range("J83").Value = 246.66667
$ws.Range("L83").Value = 2220.00003
$ws.Range("N83").Value = -12204.00003
$ws.Range("H125").Value = 1451.3334
$ws.Range("I125").Value = 765.5
$ws.Range("J125").Value = 2000
$ws.Range("K125").Value = 6889.5
$ws.Range("L125").Value = 18000
$ws.Range("M125").Value = -4429.5
$ws.Range("N125").Value = -22920
$ws.Range("H135").Value = 1741.6923
$ws.Range("I135").Value = 1741.6923
$ws.Range("K135").Value = 15675.2307
$ws.Range("M135").Value = -13140.2307
$ws.Range("H138").Value = 4959.6333
$ws.Range("J138").Value = 5116.4585
$ws.Range("L138").Value = 15349.3755
$ws.Range("N138").Value = -25629.3755

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2682.625
$ws.Range("I2").Value = 2742.75
$ws.Range("J2").Value = 2622.5
$ws.Range("K2").Value = 2742.75
$ws.Range("L2").Value = 2622.5
$ws.Range("M2").Value = -2629.75
$ws.Range("N2").Value = -2848.5
$ws.Range("H5").Value = 378.2
$ws.Range("I5").Value = 380.25
$ws.Range("K5").Value = 380.25
$ws.Range("M5").Value = -268.25
$ws.Range("H23").Value = 35997
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 35997
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 35997
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -36515
$ws.Range("H37").Value = 40088
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H45").Value = 3998
$ws.Range("I45").Value = 3997
$ws.Range("K45").Value = 3997
$ws.Range("M45").Value = -3620
$ws.Range("H63").Value = 4666.6665
$ws.Range("I63").Value = 1000
$ws.Range("J63").Value = 12000
$ws.Range("K63").Value = 1000
$ws.Range("L63").Value = 12000
$ws.Range("M63").Value = -314
$ws.Range("N63").Value = -13372
$ws.Range("H66").Value = 4666.6665
$ws.Range("I66").Value = 1000
$ws.Range("J66").Value = 12000
$ws.Range("K66").Value = 5000
$ws.Range("L66").Value = 60000
$ws.Range("M66").Value = -1568
$ws.Range("N66").Value = -66864
$ws.Range("H97").Value = 305
$ws.Range("I97").Value = 246
$ws.Range("K97").Value = 246
$ws.Range("M97").Value = 250
$ws.Range("H102").Value = 1395.6
$ws.Range("I102").Value = 1395.6
$ws.Range("K102").Value = 1395.6
$ws.Range("M102").Value = 226.4000000000001
$ws.Range("H110").Value = 1293.6666
$ws.Range("I110").Value = 1293.6666
$ws.Range("K110").Value = 1293.6666
$ws.Range("M110").Value = 751.3334
$ws.Range("H116").Value = 2682.625
$ws.Range("I116").Value = 2742.75
$ws.Range("J116").Value = 2622.5
$ws.Range("K116").Value = 2742.75
$ws.Range("L116").Value = 2622.5
$ws.Range("M116").Value = -448.75
$ws.Range("N116").Value = -7210.5
$ws.Range("H122").Value = 3533.5
$ws.Range("I122").Value = 3802.75
$ws.Range("K122").Value = 11408.25
$ws.Range("M122").Value = -8958.25
$ws.Range("H132").Value = 1515.3334
$ws.Range("I132").Value = 1660.6666
$ws.Range("K132").Value = 4981.9998
$ws.Range("M132").Value = -2451.9998

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2682.625
$ws.Range("I3").Value = 2742.75
$ws.Range("J3").Value = 2622.5
$ws.Range("K3").Value = 2742.75
$ws.Range("L3").Value = 2622.5
$ws.Range("M3").Value = -2628.75
$ws.Range("N3").Value = -2850.5
$ws.Range("H4").Value = 378.2
$ws.Range("I4").Value = 380.25
$ws.Range("K4").Value = 380.25
$ws.Range("M4").Value = -265.25
$ws.Range("H22").Value = 432.66666
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H105").Value = 7369.3
$ws.Range("I105").Value = 7336.75
$ws.Range("K105").Value = 7336.75
$ws.Range("M105").Value = -5589.75

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H33").Value = 531
$ws.Range("I33").Value = 531
$ws.Range("K33").Value = 531
$ws.Range("M33").Value = -152
$ws.Range("H62").Value = 4493.5
$ws.Range("I62").Value = 4650
$ws.Range("J62").Value = 4337
$ws.Range("K62").Value = 4650
$ws.Range("L62").Value = 4337
$ws.Range("M62").Value = -4026
$ws.Range("N62").Value = -5585
$ws.Range("H65").Value = 4493.5
$ws.Range("I65").Value = 4650
$ws.Range("J65").Value = 4337
$ws.Range("K65").Value = 23250
$ws.Range("L65").Value = 21685
$ws.Range("M65").Value = -20130
$ws.Range("N65").Value = -27925
$ws.Range("H99").Value = 9234.200000000001
$ws.Range("I99").Value = 9063.286
$ws.Range("K99").Value = 9063.286
$ws.Range("M99").Value = -7565.286
$ws.Range("H105").Value = 4624.75
$ws.Range("I105").Value = 4249.5
$ws.Range("K105").Value = 4249.5
$ws.Range("M105").Value = -2502.5
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H126").Value = 9234.200000000001
$ws.Range("I126").Value = 9063.286
$ws.Range("K126").Value = 27189.858
$ws.Range("M126").Value = -24719.858

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 279.75
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 279.75
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 839.25
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1309.25
$ws.Range("H139").Value = 1496.8572
$ws.Range("I139").Value = 1496.8572
$ws.Range("K139").Value = 4490.571599999999
$ws.Range("M139").Value = 649.4284000000007
$ws.Range("H140").Value = 1285.25
$ws.Range("I140").Value = 1285.25
$ws.Range("K140").Value = 3855.75
$ws.Range("M140").Value = 1324.25
$ws.Range("H141").Value = 6291.4
$ws.Range("I141").Value = 6291.4
$ws.Range("K141").Value = 18874.2
$ws.Range("M141").Value = -13694.2

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 71
$ws.Range("I2").Value = 71
$ws.Range("K2").Value = 71
$ws.Range("M2").Value = 42
$ws.Range("H102").Value = 2527.625
$ws.Range("I102").Value = 2367.8333
$ws.Range("K102").Value = 2367.8333
$ws.Range("M102").Value = -745.8332999999998
$ws.Range("H113").Value = 2237.8572
$ws.Range("J113").Value = 2859.8
$ws.Range("L113").Value = 2859.8
$ws.Range("N113").Value = -7199.8
$ws.Range("H132").Value = 2243.5
$ws.Range("I132").Value = 2243
$ws.Range("K132").Value = 6729
$ws.Range("M132").Value = -4199

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6093.1113
$ws.Range("I22").Value = 5223.25
$ws.Range("K22").Value = 5223.25
$ws.Range("M22").Value = -4928.25
$ws.Range("H27").Value = 6093.1113
$ws.Range("I27").Value = 5223.25
$ws.Range("K27").Value = 5223.25
$ws.Range("M27").Value = -5116.25
$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1864
$ws.Range("H122").Value = 10572.714
$ws.Range("I122").Value = 10002
$ws.Range("J122").Value = 10801
$ws.Range("K122").Value = 30006
$ws.Range("L122").Value = 32403
$ws.Range("M122").Value = -27556
$ws.Range("N122").Value = -37303

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2657
$ws.Range("I132").Value = 2079.2666
$ws.Range("K132").Value = 6237.7998
$ws.Range("M132").Value = -3707.7998
